$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change G2 from text "11.038" to the number 11
$ws.Range("G2").Value = 11

# Update the selected cell/range to G2 (was G3)
$ws.Range("G2").Select()
